# Slide 11 ("Demo") has an empty content placeholder (shape 2).
# The edit adds three blank paragraphs followed by a paragraph containing
# a YouTube link, replacing the single empty paragraph that was there.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)

$shp.TextFrame.TextRange.Text = "`r`r`rhttps://www.youtube.com/watch?v=NG29ArAEDvY"
